# updating GDP and algo
# Refresh the "Total Gross Domestic Product for San Francisco-Oakland-Hayward, CA (MSA)"
# FRED series: revised historical values (B12:B30) + one new observation (row 31, 2020-01-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised GDP values for existing years (B12:B30) ---
$ws.Range("B12").Value = 240223.33100000001
$ws.Range("B13").Value = 238126.36199999999
$ws.Range("B14").Value = 246289.538
$ws.Range("B15").Value = 256414.07
$ws.Range("B16").Value = 278166.39299999998
$ws.Range("B17").Value = 296377.34100000001
$ws.Range("B18").Value = 307604.98
$ws.Range("B19").Value = 324695.02899999998
$ws.Range("B20").Value = 309329.272
$ws.Range("B21").Value = 316779.48100000003
$ws.Range("B22").Value = 329715.79200000002
$ws.Range("B23").Value = 363995.87099999998
$ws.Range("B24").Value = 382609.21399999998
$ws.Range("B25").Value = 412207.01199999999
$ws.Range("B26").Value = 443406.28100000002
$ws.Range("B27").Value = 472395.54499999998
$ws.Range("B28").Value = 519358.95
$ws.Range("B29").Value = 562046.31099999999
$ws.Range("B30").Value = 595294.62

# --- New observation row: 2020-01-01 / 588335.543 ---
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 588335.54299999995
$ws.Range("B31").NumberFormat = "0.000"

# --- Selection state left by the editor (full A:B column selection) ---
$ws.Range("A1:B1048576").Select()
